# punto 2 tp2 WIP
# Applies the edits described by the commit: updates the "Min" solution
# texts in B10/C10/D10/C11/D11/D12, appends a small black-square "status"
# marker (rich-text run in Segoe UI) to three of those cells, updates two
# of the cell comments, tweaks row 11 height / column D width, moves the
# active selection to D12, and sets the page orientation to portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bullet = [string][char]0x25AA

function Set-BulletCell($addr, $baseText) {
    # Writes $baseText followed by a small black square "*" marker, where
    # the marker is its own rich-text run rendered in Segoe UI 11 so it
    # shows up as a distinct glyph from the rest of the (Calibri) text.
    $cell = $ws.Range($addr)
    $cell.Value = $baseText + $bullet
    $chars = $cell.Characters($baseText.Length + 1, 1)
    $chars.Font.Name = "Segoe UI"
    $chars.Font.Size = 11
}

# --- Row 10 ---------------------------------------------------------------
Set-BulletCell "B10" "Q:10, Costo=2"
$ws.Range("C10").Value = "Q:2, Costo=2+4*10+2=44"
$ws.Range("D10").Value = "Q: 9, Costo= 2+4*10+2+2=46"

# --- Row 11 -----------------------------------------------------------------
Set-BulletCell "C11" "Q:12, Costo=2+2=4"
$ws.Range("D11").Value = "Q:1, Costo=2+2+8*4+2=38"

# --- Row 12 -----------------------------------------------------------------
Set-BulletCell "D12" "Q:9, Costo=4+2=6"

# --- Comments ---------------------------------------------------------------
$commentD10 = $ws.Range("D10").Comment
$textD10 = "Maria Ines Parnisari:" + [char]10 + `
    "Comprar todo lo que se pueda en el mes 1(20), lo del mes 2(2), lo que sobre en el 3(9)" + [char]10
$null = $commentD10.Text($textD10)

$commentD11 = $ws.Range("D11").Comment
$textD11 = "Maria Ines Parnisari:" + [char]10 + `
    "Comprar lo optimo del mes 1 (10), todo lo que se pueda en el mes 2 (20), lo que sobre en el mes 3(1)"
$null = $commentD11.Text($textD11)

# --- Layout tweaks ------------------------------------------------------
# Row 11 grew slightly taller (bigger font used in one of its cells).
$ws.Rows(11).RowHeight = 16.5

# Column D widened to fit the new, longer "Min" text.
$ws.Columns(4).ColumnWidth = 24.43

# Active selection ends up on D12.
$null = $ws.Range("D12").Select()

# Page was set to portrait orientation.
$ws.PageSetup.Orientation = 1

Write-Host "applied progdinamica_pto2 edits"
